$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("4:6").Delete()
$ws.Range("U2").Value = "eyescrunching+jaw"
$ws.Range("U3").Value = "jaw+raisingeyebrows"

$ws.Range("Q2").Value = -0.02764303769192666
$ws.Range("R2").Value = 19.79815774584502
$ws.Range("S2").Value = -14.65713500675496

$ws.Range("Q3").Value = -0.02208809892450309
$ws.Range("R3").Value = 36.68747884301371
$ws.Range("S3").Value = -11.39052687350372

Write-Output "done"
